$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New transaction rows (report section) appended after the existing 3 rows.
# Columns: A = Mã giao dịch, B = Loại giao dịch, C = Ngày giao dịch, D = Số tiền, E = Ghi chú
$data = @(
    @("transt2", "Nạp tiền", "11/10/2024 12:00:00 AM", "241512413", "fsafsasd"),
    @("transt3", "Rút tiền", "12/4/2024 12:00:00 AM",  "5344223",   "sfsadfsa"),
    @("transt4", "Nạp tiền", "12/4/2024 12:00:00 AM",  "421424",    "ỉa chải"),
    @("transt5", "Rút tiền", "12/2/2024 12:00:00 AM",  "521432",    "dfasa")
)

$startRow = 4
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]

    $ws.Cells.Item($row, 1).Value = $rowData[0]
    $ws.Cells.Item($row, 2).Value = $rowData[1]
    $ws.Cells.Item($row, 3).Value = $rowData[2]

    # Số tiền looks numeric, so it must be forced to text so it keeps being
    # stored as a shared string (matching the rest of the sheet), then the
    # number-format override is cleared again so the cell keeps the default
    # style.
    $ws.Cells.Item($row, 4).NumberFormat = "@"
    $ws.Cells.Item($row, 4).Value = $rowData[3]
    $ws.Cells.Item($row, 4).ClearFormats()

    $ws.Cells.Item($row, 5).Value = $rowData[4]
}
